$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 27875.3
$ws.Range("B21").Value = 298532.49
$ws.Range("F50").Value = 167
$ws.Range("G50").Value = 32212.63
$ws.Range("F51").Value = 66
$ws.Range("G51").Value = 2779.92
$ws.Range("F54").Value = 79
$ws.Range("G54").Value = 1796.46
$ws.Range("F59").Value = 95
$ws.Range("G59").Value = 3708.8
$ws.Range("B77").Value = 127813.07
$ws.Range("F120").Value = 12
$ws.Range("G120").Value = 323879.76
$ws.Range("B122").Value = 962012.97
$ws.Range("F142").Value = 39
$ws.Range("G142").Value = 3760.77
$ws.Range("F146").Value = 69
$ws.Range("G146").Value = 1330.32
$ws.Range("B147").Value = 57324.73
$ws.Range("F152").Value = 17
$ws.Range("G152").Value = 1045.84
$ws.Range("B154").Value = 1374.19
$ws.Range("F165").Value = 29
$ws.Range("G165").Value = 2655.53
$ws.Range("B168").Value = 28949.09
$ws.Range("B179").Value = 57756
$ws.Range("B180").Value = 53925
$ws.Range("F182").Value = 2
$ws.Range("G182").Value = 209.1
$ws.Range("B183").Value = 28727.84
$ws.Range("F203").Value = 106
$ws.Range("G203").Value = 2688.16
$ws.Range("F208").Value = 29
$ws.Range("G208").Value = 1325.3
$ws.Range("F209").Value = 20
$ws.Range("G209").Value = 761.6
$ws.Range("B210").Value = 29476.82
$ws.Range("F226").Value = 9
$ws.Range("G226").Value = 464.85
$ws.Range("B237").Value = 13382.31
$ws.Range("F307").Value = 84
$ws.Range("G307").Value = 1697.64
$ws.Range("F309").Value = 161
$ws.Range("G309").Value = 5361.3
$ws.Range("F328").Value = 0
$ws.Range("G328").Value = 0
$ws.Range("B336").Value = 252625.5
$ws.Range("F339").Value = 1
$ws.Range("G339").Value = 80.02
$ws.Range("F340").Value = 11
$ws.Range("G340").Value = 810.8099999999999
$ws.Range("F359").Value = 3
$ws.Range("G359").Value = 513.99
$ws.Range("F361").Value = 153
$ws.Range("G361").Value = 6300.54
$ws.Range("B364").Value = 53335.27
$ws.Range("F392").Value = 218
$ws.Range("G392").Value = 4445.02
$ws.Range("B397").Value = 24072.5
$ws.Range("F406").Value = 42
$ws.Range("G406").Value = 6311.34
$ws.Range("B407").Value = 14176.56
$ws.Range("F412").Value = 6
$ws.Range("G412").Value = 275.82
$ws.Range("B421").Value = 29614.88
$ws.Range("F435").Value = 168
$ws.Range("G435").Value = 4131.12
$ws.Range("F436").Value = 111
$ws.Range("G436").Value = 2984.79
$ws.Range("B439").Value = 130699.02
$ws.Range("F464").Value = 629
$ws.Range("G464").Value = 8114.1
$ws.Range("B470").Value = 38901
$ws.Range("F479").Value = 789
$ws.Range("G479").Value = 10375.35
$ws.Range("F480").Value = 888
$ws.Range("G480").Value = 11375.28
$ws.Range("F482").Value = 463
$ws.Range("G482").Value = 9134.99
$ws.Range("F486").Value = 849
$ws.Range("G486").Value = 5586.42
$ws.Range("F488").Value = 274
$ws.Range("G488").Value = 5332.04
$ws.Range("F492").Value = 567
$ws.Range("G492").Value = 7456.05
$ws.Range("F493").Value = 575
$ws.Range("G493").Value = 15122.5
$ws.Range("F495").Value = 887
$ws.Range("G495").Value = 13065.51
$ws.Range("B496").Value = 164600.53
$ws.Range("F502").Value = 142
$ws.Range("G502").Value = 4120.84
$ws.Range("B515").Value = 37577.24
$ws.Range("F521").Value = 543
$ws.Range("G521").Value = 5820.96
$ws.Range("B539").Value = 115556.85
$ws.Range("F587").Value = 101
$ws.Range("G587").Value = 2634.08
$ws.Range("F597").Value = 27
$ws.Range("G597").Value = 4679.64
$ws.Range("B600").Value = 69863.83
$ws.Range("F605").Value = 49
$ws.Range("G605").Value = 1757.14
$ws.Range("F615").Value = 114
$ws.Range("G615").Value = 2725.74
$ws.Range("B625").Value = 45726.97
$ws.Range("F650").Value = 53
$ws.Range("G650").Value = 4400.06
$ws.Range("B651").Value = 241654.79
$ws.Range("F654").Value = 40
$ws.Range("G654").Value = 7120.8
$ws.Range("F655").Value = 9
$ws.Range("G655").Value = 1174.95
$ws.Range("B659").Value = 9481.65
$ws.Range("F710").Value = 30
$ws.Range("G710").Value = 1686.6
$ws.Range("B712").Value = 3697.92
$ws.Range("F721").Value = 50
$ws.Range("G721").Value = 7016.5
$ws.Range("B722").Value = 27253.31
$ws.Range("F748").Value = 160
$ws.Range("G748").Value = 21296
$ws.Range("F755").Value = 303
$ws.Range("G755").Value = 21076.68
$ws.Range("F759").Value = 170
$ws.Range("G759").Value = 25692.1
$ws.Range("B767").Value = 565819.7
$ws.Range("F783").Value = 13
$ws.Range("G783").Value = 1964.69
$ws.Range("F792").Value = 117
$ws.Range("G792").Value = 3721.77
$ws.Range("F794").Value = 54
$ws.Range("G794").Value = 2694.06
$ws.Range("B796").Value = 162841.82
$ws.Range("F848").Value = 12
$ws.Range("G848").Value = 474.36
$ws.Range("B849").Value = 1490.47
$ws.Range("B855").Value = 5303643.36
$ws.Range("B856").Value = 5303643.36
